$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in row 13 (6th log entry) that was previously blank
$ws.Range("B13").Value = 6
$ws.Range("C13").Value = 43521
$ws.Range("D13").Value = 0.79513888888888884
$ws.Range("E13").Value = 0.95833333333333337
$ws.Range("G13").Value = 165
$ws.Range("H13").Value = "Summarising"
$ws.Range("I13").Value = "MVC Notes"

# Update the active selection to H9
$ws.Range("H9").Select()
